$d = $word.ActiveDocument

function Insert-RunXml {
    param($Range, $InnerXml)
    $payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $InnerXml +
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($payload)
}

# --- Edit 1: "In the youtube video" -> mark "youtube" with spelling proof-error run split ---
$p1 = $d.Paragraphs.Item(1)
$rng1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$inner1 = '<w:p>' +
    '<w:r><w:t xml:space="preserve">In the </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>youtube</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> video</w:t></w:r>' +
    '</w:p>'
Insert-RunXml $rng1 $inner1

# --- Edit 2: "storing data in localstorage" -> mark "localstorage" as spelling proof-error ---
$p5 = $d.Paragraphs.Item(5)
$rng5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)
$inner5 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">storing data in </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>localstorage</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
Insert-RunXml $rng5 $inner5

# --- Edit 3: "checking to see if the user is logged in. If so" -> mark "checking" as grammar proof-error ---
$p6 = $d.Paragraphs.Item(6)
$rng6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$inner6 = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>checking</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> to see if the user is logged in. If so</w:t></w:r>' +
    '</w:p>'
Insert-RunXml $rng6 $inner6

# --- Edit 4: append new paragraphs after the final "When comparing strings..." paragraph, ---
# --- moving the _GoBack bookmark to the very end of the new content.                      ---
$d.Bookmarks.Item("_GoBack").Delete()

$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$endRng = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$rsquo = [string][char]0x2019
$newTail = 'The reason why I couldn' + $rsquo + 't click on an item was because the focus was on the delete button.'

$innerTail = '<w:p/>' +
    '<w:p><w:r><w:t>' + $newTail + '</w:t></w:r></w:p>' +
    '<w:p>' +
    '<w:r><w:t xml:space="preserve">Need to </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>unfocus</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> the delete button</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
Insert-RunXml $endRng $innerTail
